$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 49
$ws.Range("H49").Value = 666.6667
$ws.Range("I49").Value = 500
$ws.Range("K49").Value = 1500
$ws.Range("M49").Value = -1364

# Row 51
$ws.Range("H51").Value = 2500
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968

# Row 58
$ws.Range("H58").Value = 3132.25
$ws.Range("I58").Value = 190
$ws.Range("J58").Value = 4113
$ws.Range("K58").Value = 570
$ws.Range("L58").Value = 12339
$ws.Range("M58").Value = -420
$ws.Range("N58").Value = -12639

# Row 137
$ws.Range("H137").Value = 1046.5294
$ws.Range("I137").Value = 917.4545000000001
$ws.Range("J137").Value = 1283.1666
$ws.Range("K137").Value = 2752.3635
$ws.Range("L137").Value = 3849.4998
$ws.Range("M137").Value = -202.3635000000004
$ws.Range("N137").Value = -8949.4998

$ws = $wb.Worksheets.Item("ARM")
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Row 61
$ws.Range("H61").Value = 800.75
$ws.Range("I61").Value = 720.8
$ws.Range("K61").Value = 720.8
$ws.Range("M61").Value = -508.8

# Row 74
$ws.Range("H74").Value = 1081.5518
$ws.Range("I74").Value = 1110.2632
$ws.Range("J74").Value = 1027
$ws.Range("K74").Value = 1110.2632
$ws.Range("L74").Value = 1027
$ws.Range("M74").Value = -236.2632000000001
$ws.Range("N74").Value = -2775

# Row 77
$ws.Range("H77").Value = 1081.5518
$ws.Range("I77").Value = 1110.2632
$ws.Range("J77").Value = 1027
$ws.Range("K77").Value = 5551.316000000001
$ws.Range("L77").Value = 5135
$ws.Range("M77").Value = -1183.316000000001
$ws.Range("N77").Value = -13871

# Row 132
$ws.Range("H132").Value = 1860.75
$ws.Range("I132").Value = 828.9167
$ws.Range("J132").Value = 2892.5833
$ws.Range("K132").Value = 2486.7501
$ws.Range("L132").Value = 8677.749899999999
$ws.Range("M132").Value = 43.2498999999998
$ws.Range("N132").Value = -13737.7499

# Row 136
$ws.Range("H136").Value = 800.75
$ws.Range("I136").Value = 720.8
$ws.Range("K136").Value = 2162.4
$ws.Range("M136").Value = 387.6000000000004

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 69674.55
$ws.Range("I134").Value = 5446.476
$ws.Range("J134").Value = 204553.5
$ws.Range("K134").Value = 16339.428
$ws.Range("L134").Value = 613660.5
$ws.Range("M134").Value = -13804.428
$ws.Range("N134").Value = -618730.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1590.2941
$ws.Range("J31").Value = 2652.3333
$ws.Range("L31").Value = 2652.3333
$ws.Range("N31").Value = -3242.3333

# Row 34
$ws.Range("H34").Value = 1590.2941
$ws.Range("J34").Value = 2652.3333
$ws.Range("L34").Value = 2652.3333
$ws.Range("N34").Value = -3056.3333

# Row 58
$ws.Range("H58").Value = 4460
$ws.Range("I58").Value = 698.5217
$ws.Range("J58").Value = 21762.8
$ws.Range("K58").Value = 698.5217
$ws.Range("L58").Value = 21762.8
$ws.Range("M58").Value = -495.5217
$ws.Range("N58").Value = -22168.8

# Row 132
$ws.Range("H132").Value = 1423.963
$ws.Range("I132").Value = 879.5
$ws.Range("J132").Value = 3819.6
$ws.Range("K132").Value = 2638.5
$ws.Range("L132").Value = 11458.8
$ws.Range("M132").Value = -108.5
$ws.Range("N132").Value = -16518.8

# Row 134
$ws.Range("H134").Value = 2089.2285
$ws.Range("I134").Value = 1540.5555
$ws.Range("J134").Value = 3941
$ws.Range("K134").Value = 4621.666499999999
$ws.Range("L134").Value = 11823
$ws.Range("M134").Value = -2086.666499999999
$ws.Range("N134").Value = -16893

# Row 136
$ws.Range("H136").Value = 4460
$ws.Range("I136").Value = 698.5217
$ws.Range("J136").Value = 21762.8
$ws.Range("K136").Value = 2095.5651
$ws.Range("L136").Value = 65288.39999999999
$ws.Range("M136").Value = 454.4349000000002
$ws.Range("N136").Value = -70388.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1335
$ws.Range("I34").Value = 602
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1806
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = -1722
$ws.Range("N34").Value = -15168

# Row 39
$ws.Range("H39").Value = 2058.3684
$ws.Range("J39").Value = 2058.3684
$ws.Range("L39").Value = 6175.1052
$ws.Range("N39").Value = -6763.1052

# Row 55
$ws.Range("H55").Value = 48532.953
$ws.Range("J55").Value = 1005.4706
$ws.Range("L55").Value = 3016.4118
$ws.Range("N55").Value = -3370.4118

# Row 64
$ws.Range("H64").Value = 2035.5
$ws.Range("I64").Value = 753
$ws.Range("J64").Value = 2676.75
$ws.Range("K64").Value = 2259
$ws.Range("L64").Value = 8030.25
$ws.Range("M64").Value = -1989
$ws.Range("N64").Value = -8570.25

# Row 67
$ws.Range("H67").Value = 2035.5
$ws.Range("I67").Value = 753
$ws.Range("J67").Value = 2676.75
$ws.Range("K67").Value = 2259
$ws.Range("L67").Value = 8030.25
$ws.Range("M67").Value = -1323
$ws.Range("N67").Value = -9902.25

# Row 131
$ws.Range("H131").Value = 797.0205999999999
$ws.Range("I131").Value = 319.25
$ws.Range("J131").Value = 891.3951
$ws.Range("K131").Value = 957.75
$ws.Range("L131").Value = 2674.1853
$ws.Range("M131").Value = 4082.25
$ws.Range("N131").Value = -12754.1853

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2492.2285
$ws.Range("I132").Value = 2144.3044
$ws.Range("K132").Value = 6432.9132
$ws.Range("M132").Value = -3902.9132

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1207.7
$ws.Range("I22").Value = 1050.4286
$ws.Range("J22").Value = 1292.3846
$ws.Range("K22").Value = 1050.4286
$ws.Range("L22").Value = 1292.3846
$ws.Range("M22").Value = -755.4286
$ws.Range("N22").Value = -1882.3846

# Row 27
$ws.Range("H27").Value = 1207.7
$ws.Range("I27").Value = 1050.4286
$ws.Range("J27").Value = 1292.3846
$ws.Range("K27").Value = 1050.4286
$ws.Range("L27").Value = 1292.3846
$ws.Range("M27").Value = -943.4286
$ws.Range("N27").Value = -1506.3846

# Row 46
$ws.Range("H46").Value = 2100
$ws.Range("I46").Value = 2833.3333
$ws.Range("K46").Value = 2833.3333
$ws.Range("M46").Value = -2645.3333

# Row 132
$ws.Range("H132").Value = 3527.879
$ws.Range("I132").Value = 3191
$ws.Range("K132").Value = 9573
$ws.Range("M132").Value = -7043

# Row 136
$ws.Range("H136").Value = 3950.8708
$ws.Range("I136").Value = 1429.2174
$ws.Range("K136").Value = 4287.6522
$ws.Range("M136").Value = -1737.6522

$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 14800
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 14800
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 14800
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -15416

# Row 132
$ws.Range("H132").Value = 2681.5386
$ws.Range("I132").Value = 2587.6
$ws.Range("J132").Value = 2994.6667
$ws.Range("K132").Value = 7762.799999999999
$ws.Range("L132").Value = 8984.000100000001
$ws.Range("M132").Value = -5232.799999999999
$ws.Range("N132").Value = -14044.0001

# Row 136
$ws.Range("H136").Value = 1598.5294
$ws.Range("I136").Value = 1629.4828
$ws.Range("K136").Value = 4888.4484
$ws.Range("M136").Value = -2338.4484
